$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.817.47"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.031.11"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "2.332.32"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "2.044.51"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "37.780.20"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.128"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0601"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.41%  "
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "1.536.96"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0919"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.13%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "2.222.56"
$ws.Range("E51").Value = "  -1.30%  "
